# Update the "想去人数" (interest count, column F) figures on the
# "展览" and "全部类型" sheets to reflect the latest generated data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 573
$ws1.Range("F3").Value = 190
$ws1.Range("F4").Value = 376
$ws1.Range("F5").Value = 424
$ws1.Range("F6").Value = 266
$ws1.Range("F7").Value = 2418
$ws1.Range("F9").Value = 6326
$ws1.Range("F11").Value = 408

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 573
$ws4.Range("F3").Value = 190
$ws4.Range("F4").Value = 376
$ws4.Range("F5").Value = 424
$ws4.Range("F6").Value = 266
$ws4.Range("F9").Value = 2418
$ws4.Range("F11").Value = 6326
$ws4.Range("F13").Value = 408
